$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2
$ws.Range("D2").Value = "27.897.65"
$ws.Range("E2").Value = "  +1.37%  "

# Row 3
$ws.Range("D3").Value = "1.760.41"
$ws.Range("E3").Value = "  +0.17%  "

# Row 4
Set-TextValue $ws.Range("D4") "1.000"
$ws.Range("E4").Value = "  -0.49%  "

# Row 5
Set-TextValue $ws.Range("D5") "322.12"
$ws.Range("E5").Value = "  -0.83%  "

# Row 6
Set-TextValue $ws.Range("D6") "0.9976"
$ws.Range("E6").Value = "  -0.54%  "

# Row 7
Set-TextValue $ws.Range("D7") "0.4255"
$ws.Range("E7").Value = "  -4.00%  "

# Row 8
Set-TextValue $ws.Range("D8") "0.3599"
$ws.Range("E8").Value = "  -2.49%  "

# Row 9
$ws.Range("E9").Value = "  -1.15%  "

# Row 10
Set-TextValue $ws.Range("D10") "0.07443"
$ws.Range("E10").Value = "  -2.77%  "

# Row 11
Set-TextValue $ws.Range("D11") "1.101"
$ws.Range("E11").Value = "  -0.92%  "

# Row 12
Set-TextValue $ws.Range("D12") "0.9990"
$ws.Range("E12").Value = "  -0.37%  "

# Row 13
Set-TextValue $ws.Range("D13") "21.47"
$ws.Range("E13").Value = "  -0.56%  "

# Row 14
$ws.Range("E14").Value = "  -0.70%  "

# Row 15
$ws.Range("E15").Value = "  -1.89%  "

# Row 16
$ws.Range("D16").Value = "1.794.30"
$ws.Range("E16").Value = "  +1.99%  "

# Row 17
Set-TextValue $ws.Range("D17") "91.09"
$ws.Range("E17").Value = "  +0.98%  "

# Row 18
Set-TextValue $ws.Range("D18") "0.00001059"
$ws.Range("E18").Value = "  -1.09%  "

# Row 19
Set-TextValue $ws.Range("D19") "0.06375"
$ws.Range("E19").Value = "  +1.84%  "

# Row 20
Set-TextValue $ws.Range("D20") "0.9996"
$ws.Range("E20").Value = "  -0.28%  "

# Row 21
Set-TextValue $ws.Range("D21") "17.15"
$ws.Range("E21").Value = "  -1.22%  "

# Row 22
Set-TextValue $ws.Range("D22") "5.957"
$ws.Range("E22").Value = "  -3.38%  "

# Row 23
$ws.Range("D23").Value = "27.913.64"
$ws.Range("E23").Value = "  +1.33%  "

# Row 24
Set-TextValue $ws.Range("D24") "11.31"
$ws.Range("E24").Value = "  -1.94%  "

# Row 25
Set-TextValue $ws.Range("D25") "2.137"
$ws.Range("E25").Value = "  -7.38%  "

# Row 26
Set-TextValue $ws.Range("D26") "160.68"
$ws.Range("E26").Value = "  +5.15%  "

# Row 27
Set-TextValue $ws.Range("D27") "20.26"
$ws.Range("E27").Value = "  -1.28%  "

# Row 28
$ws.Range("D28").Value = "1.992.59"
$ws.Range("E28").Value = "  +1.84%  "

# Row 29
Set-TextValue $ws.Range("D29") "2.143"
$ws.Range("E29").Value = "  -6.34%  "

# Row 30
Set-TextValue $ws.Range("D30") "126.04"
$ws.Range("E30").Value = "  -1.24%  "

# Row 31
Set-TextValue $ws.Range("D31") "1.171"
$ws.Range("E31").Value = "  -0.26%  "

# Row 32
Set-TextValue $ws.Range("D32") "5.676"
$ws.Range("E32").Value = "  -0.56%  "

# Row 33
Set-TextValue $ws.Range("D33") "0.09011"
$ws.Range("E33").Value = "  -2.01%  "

# Row 34
Set-TextValue $ws.Range("D34") "3.511"
$ws.Range("E34").Value = "  -3.05%  "

# Row 35
Set-TextValue $ws.Range("D35") "12.65"
$ws.Range("E35").Value = "  +0.61%  "

# Row 36
Set-TextValue $ws.Range("D36") "0.02309"
$ws.Range("E36").Value = "  +0.05%  "

# Row 37
Set-TextValue $ws.Range("D37") "0.06085"
$ws.Range("E37").Value = "  -0.12%  "

# Row 38
$ws.Range("B38").Value = "InternetComputer(DFINITY)"
$ws.Range("C38").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue $ws.Range("D38") "5.048"
$ws.Range("E38").Value = "  +0.11%  "

# Row 39
$ws.Range("B39").Value = "Algorand"
$ws.Range("C39").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue $ws.Range("D39") "0.2110"
$ws.Range("E39").Value = "  -2.05%  "

# Row 40
Set-TextValue $ws.Range("D40") "0.6402"
$ws.Range("E40").Value = "  -0.41%  "

# Row 41
Set-TextValue $ws.Range("D41") "1.183"
$ws.Range("E41").Value = "  +0.77%  "

# Row 42
Set-TextValue $ws.Range("D42") "1.000"
$ws.Range("E42").Value = "  -0.24%  "

# Row 43
Set-TextValue $ws.Range("D43") "7.896"
$ws.Range("E43").Value = "  -0.65%  "

# Row 44
Set-TextValue $ws.Range("D44") "1.392"
$ws.Range("E44").Value = "  -1.24%  "

# Row 45
Set-TextValue $ws.Range("D45") "13.64"
$ws.Range("E45").Value = "  -0.47%  "

# Row 46
Set-TextValue $ws.Range("D46") "0.5960"
$ws.Range("E46").Value = "  +0.09%  "

# Row 47
Set-TextValue $ws.Range("D47") "3.702"
$ws.Range("E47").Value = "  -0.52%  "

# Row 48
Set-TextValue $ws.Range("D48") "124.45"
$ws.Range("E48").Value = "  -1.23%  "

# Row 49
Set-TextValue $ws.Range("D49") "1.983"
$ws.Range("E49").Value = "  +0.22%  "

# Row 50
$ws.Range("E50").Value = "  +1.26%  "

# Row 51
$ws.Range("E51").Value = "  -0.06%  "
